$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New journal entry on row 28 (continuing the Tableau1 table) ---

# Date (col B) - copy formatting (short date style) from the row above, then set the value
$ws.Range("B27").Copy()
$ws.Range("B28").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B28").Value = 44266

# Heure début / Heure fin (cols C/D) - column already carries the time style
$ws.Range("C28").Value = 0.375
$ws.Range("D28").Value = 0.38541666666666669

# Durée (col E) - recreate the calculated column formula for this row
$ws.Range("E28").Formula = '=IF(ISBLANK(Tableau1[[#This Row],[Heure fin]]),"",Tableau1[[#This Row],[Heure fin]]-Tableau1[[#This Row],[Heure début]])'

# Module / Type / Tâche / Lieu / Descriptif / Terminer
$ws.Range("F28").Value = "Ma-20"
$ws.Range("G28").Value = "Code"
$ws.Range("H28").Value = "Jeu"
$ws.Range("I28").Value = "CPNV"

$ws.Range("J27").Copy()
$ws.Range("J28").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("J28").Value = "J'ai fait que les casse tocuher avec un bateau soit différente qu'une case touchée avec rien"

$ws.Range("K28").Value = "Oui"

# Row grew to fit the wrapped description text
$ws.Rows.Item(28).RowHeight = 43.2

# Selection left where the author ended up working
$null = $ws.Range("G29").Select()

$wb.Save()
